$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "D2 before:" $ws.Range("D2").Value
$ws.Range("D2").Value = "Capstone Project "
Write-Host "D2 after:" $ws.Range("D2").Value
Write-Host "C8 before:" $ws.Range("C8").Value
$ws.Range("C8").Value = "Project Conception & Inititation "
Write-Host "E8 set date"
$ws.Range("E8").Value = "8/30/2017"
Write-Host "E8 after:" $ws.Range("E8").Value
